$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 33.65726433333334
$ws.Range("H2").Value = 100.971793
$ws.Range("I2").Value = 0.8115737688004754
$ws.Range("J2").Value = 0.8115737688004754
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.83081766666667
$ws.Range("N2").Value = 80.49245300000001
$ws.Range("O2").Value = 0.5916656861001716
$ws.Range("P2").Value = 0.5916656861001716
$ws.Range("Q2").Value = 903.0519224864702
$ws.Range("R2").Value = 8127.467302378231
$ws.Range("S2").Value = 0.4801803507382353
$ws.Range("T2").Value = 0.4801803507382353

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 33.65726433333334
$ws.Range("H3").Value = 100.971793
$ws.Range("I3").Value = 0.8115737688004754
$ws.Range("J3").Value = 0.8115737688004754
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.21969166666667
$ws.Range("N3").Value = 30.659075
$ws.Range("O3").Value = 0.2253617819930474
$ws.Range("P3").Value = 0.2253617819930474
$ws.Range("Q3").Value = 343.966863830164
$ws.Range("R3").Value = 3095.701774471475
$ws.Range("S3").Value = 0.1828977107556886
$ws.Range("T3").Value = 0.1828977107556886

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 33.65726433333334
$ws.Range("H4").Value = 100.971793
$ws.Range("I4").Value = 0.8115737688004754
$ws.Range("J4").Value = 0.8115737688004754
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.297426666666667
$ws.Range("N4").Value = 24.89228
$ws.Range("O4").Value = 0.1829725319067811
$ws.Range("P4").Value = 0.1829725319067811
$ws.Range("Q4").Value = 279.2686826064489
$ws.Range("R4").Value = 2513.41814345804
$ws.Range("S4").Value = 0.1484957073065515
$ws.Range("T4").Value = 0.1484957073065515

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.56955
$ws.Range("H5").Value = 1.70865
$ws.Range("I5").Value = 0.01373349406661455
$ws.Range("J5").Value = 0.01373349406661455
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.83081766666667
$ws.Range("N5").Value = 80.49245300000001
$ws.Range("O5").Value = 0.5916656861001716
$ws.Range("P5").Value = 0.5916656861001716
$ws.Range("Q5").Value = 15.28149220205
$ws.Range("R5").Value = 137.53342981845
$ws.Range("S5").Value = 0.00812563718947613
$ws.Range("T5").Value = 0.00812563718947613

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.56955
$ws.Range("H6").Value = 1.70865
$ws.Range("I6").Value = 0.01373349406661455
$ws.Range("J6").Value = 0.01373349406661455
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.21969166666667
$ws.Range("N6").Value = 30.659075
$ws.Range("O6").Value = 0.2253617819930474
$ws.Range("P6").Value = 0.2253617819930474
$ws.Range("Q6").Value = 5.820625388750001
$ws.Range("R6").Value = 52.38562849875
$ws.Range("S6").Value = 0.003095004695843197
$ws.Range("T6").Value = 0.003095004695843197

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.56955
$ws.Range("H7").Value = 1.70865
$ws.Range("I7").Value = 0.01373349406661455
$ws.Range("J7").Value = 0.01373349406661455
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.297426666666667
$ws.Range("N7").Value = 24.89228
$ws.Range("O7").Value = 0.1829725319067811
$ws.Range("P7").Value = 0.1829725319067811
$ws.Range("Q7").Value = 4.725799358
$ws.Range("R7").Value = 42.532194222
$ws.Range("S7").Value = 0.002512852181295218
$ws.Range("T7").Value = 0.002512852181295218

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.244787666666667
$ws.Range("H8").Value = 21.734363
$ws.Range("I8").Value = 0.1746927371329101
$ws.Range("J8").Value = 0.1746927371329101
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.83081766666667
$ws.Range("N8").Value = 80.49245300000001
$ws.Range("O8").Value = 0.5916656861001716
$ws.Range("P8").Value = 0.5916656861001716
$ws.Range("Q8").Value = 194.3835769180488
$ws.Range("R8").Value = 1749.452192262439
$ws.Range("S8").Value = 0.1033596981724601
$ws.Range("T8").Value = 0.1033596981724601

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.244787666666667
$ws.Range("H9").Value = 21.734363
$ws.Range("I9").Value = 0.1746927371329101
$ws.Range("J9").Value = 0.1746927371329101
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.21969166666667
$ws.Range("N9").Value = 30.659075
$ws.Range("O9").Value = 0.2253617819930474
$ws.Range("P9").Value = 0.2253617819930474
$ws.Range("Q9").Value = 74.0394961438028
$ws.Range("R9").Value = 666.355465294225
$ws.Range("S9").Value = 0.03936906654151561
$ws.Range("T9").Value = 0.03936906654151561

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.244787666666667
$ws.Range("H10").Value = 21.734363
$ws.Range("I10").Value = 0.1746927371329101
$ws.Range("J10").Value = 0.1746927371329101
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.297426666666667
$ws.Range("N10").Value = 24.89228
$ws.Range("O10").Value = 0.1829725319067811
$ws.Range("P10").Value = 0.1829725319067811
$ws.Range("Q10").Value = 60.11309437973778
$ws.Range("R10").Value = 541.01784941764
$ws.Range("S10").Value = 0.0319639724189343
$ws.Range("T10").Value = 0.0319639724189343
